$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (column C) date from 2023-11-21 (45251) to 2023-11-22 (45252)
# for every existing data row (rows 2 through 26).
for ($r = 2; $r -le 26; $r++) {
    $ws.Cells.Item($r, 3).Value = 45252
}

# Row 26 picks up an explicit row height (matching the sheet default of 15).
$ws.Rows.Item(26).RowHeight = 15

# Append new row 27: A 58731-2023
$ws.Range("A27").Value = "A 58731-2023"
$ws.Range("B27").Value = 45251
$ws.Range("B27").NumberFormat = "YYYY-MM-DD"
$ws.Range("C27").Value = 45252
$ws.Range("C27").NumberFormat = "YYYY-MM-DD"
$ws.Range("D27").Value = "OKÄNT"
$ws.Range("E27").Value = "OKÄNT"
$ws.Range("G27").Value = 10.6
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = 0
$ws.Range("N27").Value = 0
$ws.Range("O27").Value = 0
$ws.Range("P27").Value = 0
$ws.Range("Q27").Value = 0
$ws.Range("R27").Value = ""
$ws.Range("R27").WrapText = $true
$ws.Rows.Item(27).RowHeight = 15

# Append new row 28: A 58730-2023 (keeps the sheet's default row height)
$ws.Range("A28").Value = "A 58730-2023"
$ws.Range("B28").Value = 45251
$ws.Range("B28").NumberFormat = "YYYY-MM-DD"
$ws.Range("C28").Value = 45252
$ws.Range("C28").NumberFormat = "YYYY-MM-DD"
$ws.Range("D28").Value = "OKÄNT"
$ws.Range("E28").Value = "OKÄNT"
$ws.Range("G28").Value = 0.9
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 0
$ws.Range("N28").Value = 0
$ws.Range("O28").Value = 0
$ws.Range("P28").Value = 0
$ws.Range("Q28").Value = 0
$ws.Range("R28").Value = ""
$ws.Range("R28").WrapText = $true
